# Update recomputed correlation matrix values in the macro section block
# (rows/cols iqrMed..varMean, i.e. spreadsheet columns D:M, and their symmetric counterparts)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "-0.1609705267258305" -as [double]
$ws.Range("E2").Value = "-0.2146382865242739" -as [double]
$ws.Range("F2").Value = "-0.2390098332268563" -as [double]
$ws.Range("G2").Value = "-0.1318081117631396" -as [double]
$ws.Range("H2").Value = "-0.1751586568075148" -as [double]
$ws.Range("I2").Value = "0.02463671945219919" -as [double]
$ws.Range("J2").Value = "-0.2442575912937952" -as [double]
$ws.Range("K2").Value = "-0.2755229145892673" -as [double]
$ws.Range("L2").Value = "0.009011369349083607" -as [double]
$ws.Range("M2").Value = "0.0925382074692313" -as [double]
# Row 3
$ws.Range("D3").Value = "0.09593726714713641" -as [double]
$ws.Range("E3").Value = "0.006695480955382251" -as [double]
$ws.Range("F3").Value = "0.1790023632918195" -as [double]
$ws.Range("G3").Value = "0.236592323055357" -as [double]
$ws.Range("H3").Value = "0.1134940197879655" -as [double]
$ws.Range("I3").Value = "-0.03114961804922891" -as [double]
$ws.Range("J3").Value = "-0.0457566048628302" -as [double]
$ws.Range("K3").Value = "0.142999274089144" -as [double]
$ws.Range("L3").Value = "0.1737101849475536" -as [double]
$ws.Range("M3").Value = "-0.01979505732831893" -as [double]
# Row 4
$ws.Range("B4").Value = "-0.1609705267258305" -as [double]
$ws.Range("C4").Value = "0.09593726714713641" -as [double]
$ws.Range("E4").Value = "0.6468349212578923" -as [double]
$ws.Range("F4").Value = "0.1517513741527407" -as [double]
$ws.Range("G4").Value = "0.8005362249079484" -as [double]
$ws.Range("H4").Value = "0.992481873816037" -as [double]
$ws.Range("I4").Value = "0.879975068712724" -as [double]
$ws.Range("J4").Value = "0.7965766664278233" -as [double]
$ws.Range("K4").Value = "0.247022406716849" -as [double]
$ws.Range("L4").Value = "0.6139816498137851" -as [double]
$ws.Range("M4").Value = "0.7794274091831108" -as [double]
$ws.Range("N4").Value = "-5.26312476638483e-16" -as [double]
$ws.Range("O4").Value = "-0.4844866619882676" -as [double]
$ws.Range("R4").Value = "-0.03589309818135301" -as [double]
$ws.Range("S4").Value = "-0.06869523931865897" -as [double]
$ws.Range("T4").Value = "0.4051407095347814" -as [double]
$ws.Range("U4").Value = "-0.3079389633548364" -as [double]
# Row 5
$ws.Range("B5").Value = "-0.2146382865242739" -as [double]
$ws.Range("C5").Value = "0.006695480955382251" -as [double]
$ws.Range("D5").Value = "0.6468349212578923" -as [double]
$ws.Range("F5").Value = "0.5872880800836626" -as [double]
$ws.Range("G5").Value = "0.2893474286312728" -as [double]
$ws.Range("H5").Value = "0.6672795322378595" -as [double]
$ws.Range("I5").Value = "0.5512808457963762" -as [double]
$ws.Range("J5").Value = "0.9138020468399077" -as [double]
$ws.Range("K5").Value = "0.5629301671802283" -as [double]
$ws.Range("L5").Value = "0.2737963524037204" -as [double]
$ws.Range("M5").Value = "0.5140728904674001" -as [double]
$ws.Range("N5").Value = "-7.810849885615325e-15" -as [double]
$ws.Range("O5").Value = "-0.2130141436281778" -as [double]
$ws.Range("R5").Value = "-0.280230938060934" -as [double]
$ws.Range("S5").Value = "0.5469064887048238" -as [double]
$ws.Range("T5").Value = "0.4333987620781303" -as [double]
$ws.Range("U5").Value = "0.3026222090640849" -as [double]
# Row 6
$ws.Range("B6").Value = "-0.2390098332268563" -as [double]
$ws.Range("C6").Value = "0.1790023632918195" -as [double]
$ws.Range("D6").Value = "0.1517513741527407" -as [double]
$ws.Range("E6").Value = "0.5872880800836626" -as [double]
$ws.Range("G6").Value = "-0.09770875581285839" -as [double]
$ws.Range("H6").Value = "0.1639019505469678" -as [double]
$ws.Range("I6").Value = "0.0003716899051861331" -as [double]
$ws.Range("J6").Value = "0.502855144844632" -as [double]
$ws.Range("K6").Value = "0.9588684709491565" -as [double]
$ws.Range("L6").Value = "-0.1493617584562276" -as [double]
$ws.Range("M6").Value = "-0.006990666893042056" -as [double]
$ws.Range("N6").Value = "2.57649465002784e-16" -as [double]
$ws.Range("O6").Value = "0.0855701284130547" -as [double]
$ws.Range("R6").Value = "-0.06116014632225203" -as [double]
$ws.Range("S6").Value = "0.6117187837479278" -as [double]
$ws.Range("T6").Value = "0.2210368374975278" -as [double]
$ws.Range("U6").Value = "0.3544617503524454" -as [double]
# Row 7
$ws.Range("B7").Value = "-0.1318081117631396" -as [double]
$ws.Range("C7").Value = "0.236592323055357" -as [double]
$ws.Range("D7").Value = "0.8005362249079484" -as [double]
$ws.Range("E7").Value = "0.2893474286312728" -as [double]
$ws.Range("F7").Value = "-0.09770875581285839" -as [double]
$ws.Range("H7").Value = "0.7964976254513817" -as [double]
$ws.Range("I7").Value = "0.8108052446228973" -as [double]
$ws.Range("J7").Value = "0.4938463860778851" -as [double]
$ws.Range("K7").Value = "-0.02856147293314673" -as [double]
$ws.Range("L7").Value = "0.7920023106189256" -as [double]
$ws.Range("M7").Value = "0.719852674429839" -as [double]
$ws.Range("N7").Value = "9.535275711450758e-16" -as [double]
$ws.Range("O7").Value = "-0.5470942344968723" -as [double]
$ws.Range("R7").Value = "0.1608069006810534" -as [double]
$ws.Range("S7").Value = "-0.4059303727331163" -as [double]
$ws.Range("T7").Value = "0.1992734032194558" -as [double]
$ws.Range("U7").Value = "-0.5436623410446169" -as [double]
# Row 8
$ws.Range("B8").Value = "-0.1751586568075148" -as [double]
$ws.Range("C8").Value = "0.1134940197879655" -as [double]
$ws.Range("D8").Value = "0.992481873816037" -as [double]
$ws.Range("E8").Value = "0.6672795322378595" -as [double]
$ws.Range("F8").Value = "0.1639019505469678" -as [double]
$ws.Range("G8").Value = "0.7964976254513817" -as [double]
$ws.Range("I8").Value = "0.8831410516690518" -as [double]
$ws.Range("J8").Value = "0.8062050703183056" -as [double]
$ws.Range("K8").Value = "0.2679550097626913" -as [double]
$ws.Range("L8").Value = "0.6081518477431341" -as [double]
$ws.Range("M8").Value = "0.7908801559917824" -as [double]
$ws.Range("N8").Value = "3.144896282365325e-16" -as [double]
$ws.Range("O8").Value = "-0.468465924749046" -as [double]
$ws.Range("R8").Value = "-0.03594240137942541" -as [double]
$ws.Range("S8").Value = "-0.06505595629729606" -as [double]
$ws.Range("T8").Value = "0.3933974245641507" -as [double]
$ws.Range("U8").Value = "-0.3064871945736838" -as [double]
# Row 9
$ws.Range("B9").Value = "0.02463671945219919" -as [double]
$ws.Range("C9").Value = "-0.03114961804922891" -as [double]
$ws.Range("D9").Value = "0.879975068712724" -as [double]
$ws.Range("E9").Value = "0.5512808457963762" -as [double]
$ws.Range("F9").Value = "0.0003716899051861331" -as [double]
$ws.Range("G9").Value = "0.8108052446228973" -as [double]
$ws.Range("H9").Value = "0.8831410516690518" -as [double]
$ws.Range("J9").Value = "0.7085765328289897" -as [double]
$ws.Range("K9").Value = "0.07508672573744717" -as [double]
$ws.Range("L9").Value = "0.7775634838244121" -as [double]
$ws.Range("M9").Value = "0.9678760915982535" -as [double]
$ws.Range("N9").Value = "8.12849579998413e-15" -as [double]
$ws.Range("O9").Value = "-0.5546927522474058" -as [double]
$ws.Range("R9").Value = "-0.2162934290481094" -as [double]
$ws.Range("S9").Value = "-0.1304860423842069" -as [double]
$ws.Range("T9").Value = "0.300194367852684" -as [double]
$ws.Range("U9").Value = "-0.306152546074295" -as [double]
# Row 10
$ws.Range("B10").Value = "-0.2442575912937952" -as [double]
$ws.Range("C10").Value = "-0.0457566048628302" -as [double]
$ws.Range("D10").Value = "0.7965766664278233" -as [double]
$ws.Range("E10").Value = "0.9138020468399077" -as [double]
$ws.Range("F10").Value = "0.502855144844632" -as [double]
$ws.Range("G10").Value = "0.4938463860778851" -as [double]
$ws.Range("H10").Value = "0.8062050703183056" -as [double]
$ws.Range("I10").Value = "0.7085765328289897" -as [double]
$ws.Range("K10").Value = "0.55554696291363" -as [double]
$ws.Range("L10").Value = "0.3727458451678964" -as [double]
$ws.Range("M10").Value = "0.6266383697311377" -as [double]
$ws.Range("N10").Value = "5.597603306782835e-15" -as [double]
$ws.Range("O10").Value = "-0.2970439799879107" -as [double]
$ws.Range("R10").Value = "-0.2043475672720929" -as [double]
$ws.Range("S10").Value = "0.4155624733203104" -as [double]
$ws.Range("T10").Value = "0.4637893251987291" -as [double]
$ws.Range("U10").Value = "0.1648409968246075" -as [double]
# Row 11
$ws.Range("B11").Value = "-0.2755229145892673" -as [double]
$ws.Range("C11").Value = "0.142999274089144" -as [double]
$ws.Range("D11").Value = "0.247022406716849" -as [double]
$ws.Range("E11").Value = "0.5629301671802283" -as [double]
$ws.Range("F11").Value = "0.9588684709491565" -as [double]
$ws.Range("G11").Value = "-0.02856147293314673" -as [double]
$ws.Range("H11").Value = "0.2679550097626913" -as [double]
$ws.Range("I11").Value = "0.07508672573744717" -as [double]
$ws.Range("J11").Value = "0.55554696291363" -as [double]
$ws.Range("L11").Value = "-0.1697350189515449" -as [double]
$ws.Range("M11").Value = "0.04463684494335034" -as [double]
$ws.Range("N11").Value = "-3.663853772685152e-17" -as [double]
$ws.Range("O11").Value = "0.1048326526071779" -as [double]
$ws.Range("R11").Value = "-0.02022194646378487" -as [double]
$ws.Range("S11").Value = "0.5449871295609625" -as [double]
$ws.Range("T11").Value = "0.2412286367885081" -as [double]
$ws.Range("U11").Value = "0.2638080676402369" -as [double]
# Row 12
$ws.Range("B12").Value = "0.009011369349083607" -as [double]
$ws.Range("C12").Value = "0.1737101849475536" -as [double]
$ws.Range("D12").Value = "0.6139816498137851" -as [double]
$ws.Range("E12").Value = "0.2737963524037204" -as [double]
$ws.Range("F12").Value = "-0.1493617584562276" -as [double]
$ws.Range("G12").Value = "0.7920023106189256" -as [double]
$ws.Range("H12").Value = "0.6081518477431341" -as [double]
$ws.Range("I12").Value = "0.7775634838244121" -as [double]
$ws.Range("J12").Value = "0.3727458451678964" -as [double]
$ws.Range("K12").Value = "-0.1697350189515449" -as [double]
$ws.Range("M12").Value = "0.7920607960994264" -as [double]
$ws.Range("N12").Value = "-6.417188599869865e-15" -as [double]
$ws.Range("O12").Value = "-0.4530916816463709" -as [double]
$ws.Range("R12").Value = "-0.1030604501630649" -as [double]
$ws.Range("S12").Value = "-0.2610189723011284" -as [double]
$ws.Range("T12").Value = "0.1513992691542027" -as [double]
$ws.Range("U12").Value = "-0.3676193604901473" -as [double]
# Row 13
$ws.Range("B13").Value = "0.0925382074692313" -as [double]
$ws.Range("C13").Value = "-0.01979505732831893" -as [double]
$ws.Range("D13").Value = "0.7794274091831108" -as [double]
$ws.Range("E13").Value = "0.5140728904674001" -as [double]
$ws.Range("F13").Value = "-0.006990666893042056" -as [double]
$ws.Range("G13").Value = "0.719852674429839" -as [double]
$ws.Range("H13").Value = "0.7908801559917824" -as [double]
$ws.Range("I13").Value = "0.9678760915982535" -as [double]
$ws.Range("J13").Value = "0.6266383697311377" -as [double]
$ws.Range("K13").Value = "0.04463684494335034" -as [double]
$ws.Range("L13").Value = "0.7920607960994264" -as [double]
$ws.Range("N13").Value = "-1.202197081870329e-16" -as [double]
$ws.Range("O13").Value = "-0.5644056674894997" -as [double]
$ws.Range("R13").Value = "-0.2804808943433733" -as [double]
$ws.Range("S13").Value = "-0.115748703784505" -as [double]
$ws.Range("T13").Value = "0.2140391375821556" -as [double]
$ws.Range("U13").Value = "-0.261741111092062" -as [double]
# Row 14
$ws.Range("D14").Value = "-5.26312476638483e-16" -as [double]
$ws.Range("E14").Value = "-7.810849885615325e-15" -as [double]
$ws.Range("F14").Value = "2.57649465002784e-16" -as [double]
$ws.Range("G14").Value = "9.535275711450758e-16" -as [double]
$ws.Range("H14").Value = "3.144896282365325e-16" -as [double]
$ws.Range("I14").Value = "8.12849579998413e-15" -as [double]
$ws.Range("J14").Value = "5.597603306782835e-15" -as [double]
$ws.Range("K14").Value = "-3.663853772685152e-17" -as [double]
$ws.Range("L14").Value = "-6.417188599869865e-15" -as [double]
$ws.Range("M14").Value = "-1.202197081870329e-16" -as [double]
# Row 15
$ws.Range("D15").Value = "-0.4844866619882676" -as [double]
$ws.Range("E15").Value = "-0.2130141436281778" -as [double]
$ws.Range("F15").Value = "0.0855701284130547" -as [double]
$ws.Range("G15").Value = "-0.5470942344968723" -as [double]
$ws.Range("H15").Value = "-0.468465924749046" -as [double]
$ws.Range("I15").Value = "-0.5546927522474058" -as [double]
$ws.Range("J15").Value = "-0.2970439799879107" -as [double]
$ws.Range("K15").Value = "0.1048326526071779" -as [double]
$ws.Range("L15").Value = "-0.4530916816463709" -as [double]
$ws.Range("M15").Value = "-0.5644056674894997" -as [double]
# Row 18
$ws.Range("D18").Value = "-0.03589309818135301" -as [double]
$ws.Range("E18").Value = "-0.280230938060934" -as [double]
$ws.Range("F18").Value = "-0.06116014632225203" -as [double]
$ws.Range("G18").Value = "0.1608069006810534" -as [double]
$ws.Range("H18").Value = "-0.03594240137942541" -as [double]
$ws.Range("I18").Value = "-0.2162934290481094" -as [double]
$ws.Range("J18").Value = "-0.2043475672720929" -as [double]
$ws.Range("K18").Value = "-0.02022194646378487" -as [double]
$ws.Range("L18").Value = "-0.1030604501630649" -as [double]
$ws.Range("M18").Value = "-0.2804808943433733" -as [double]
# Row 19
$ws.Range("D19").Value = "-0.06869523931865897" -as [double]
$ws.Range("E19").Value = "0.5469064887048238" -as [double]
$ws.Range("F19").Value = "0.6117187837479278" -as [double]
$ws.Range("G19").Value = "-0.4059303727331163" -as [double]
$ws.Range("H19").Value = "-0.06505595629729606" -as [double]
$ws.Range("I19").Value = "-0.1304860423842069" -as [double]
$ws.Range("J19").Value = "0.4155624733203104" -as [double]
$ws.Range("K19").Value = "0.5449871295609625" -as [double]
$ws.Range("L19").Value = "-0.2610189723011284" -as [double]
$ws.Range("M19").Value = "-0.115748703784505" -as [double]
# Row 20
$ws.Range("D20").Value = "0.4051407095347814" -as [double]
$ws.Range("E20").Value = "0.4333987620781303" -as [double]
$ws.Range("F20").Value = "0.2210368374975278" -as [double]
$ws.Range("G20").Value = "0.1992734032194558" -as [double]
$ws.Range("H20").Value = "0.3933974245641507" -as [double]
$ws.Range("I20").Value = "0.300194367852684" -as [double]
$ws.Range("J20").Value = "0.4637893251987291" -as [double]
$ws.Range("K20").Value = "0.2412286367885081" -as [double]
$ws.Range("L20").Value = "0.1513992691542027" -as [double]
$ws.Range("M20").Value = "0.2140391375821556" -as [double]
# Row 21
$ws.Range("D21").Value = "-0.3079389633548364" -as [double]
$ws.Range("E21").Value = "0.3026222090640849" -as [double]
$ws.Range("F21").Value = "0.3544617503524454" -as [double]
$ws.Range("G21").Value = "-0.5436623410446169" -as [double]
$ws.Range("H21").Value = "-0.3064871945736838" -as [double]
$ws.Range("I21").Value = "-0.306152546074295" -as [double]
$ws.Range("J21").Value = "0.1648409968246075" -as [double]
$ws.Range("K21").Value = "0.2638080676402369" -as [double]
$ws.Range("L21").Value = "-0.3676193604901473" -as [double]
$ws.Range("M21").Value = "-0.261741111092062" -as [double]
